# "intégration des cartes du block les plus populaires"
#
# This script fills in the journal de bord entries (rows 6-10) on the
# "étape 2" sheet for the HOSTING / most_pop_hosting blocks, adjusts the
# column A width to fit the new (longer) task labels, and updates the
# saved cell-selection state on both sheets to match where the user's
# cursor ended up.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "étape 1"
$ws2 = $wb.Worksheets.Item(2)   # "étape 2"

# --- Row 6: intégration du titre du block most_pop_hosting ---
$ws2.Range("A6").Value = "intégration du titre du block most_pop_hosting"
$ws2.Range("B6").Value = "mission"
$ws2.Range("C6").Value = "pratique délibérée"
$ws2.Range("D6").Value = 44259
$ws2.Range("E6").Value = 44259
$ws2.Range("F6").Value = 1
$ws2.Range("H6").Value = 1

# --- Row 7: intégration des cartes du block most_pop_hosting ---
$ws2.Range("A7").Value = "intégration des cartes du block most_pop_hosting"
$ws2.Range("B7").Value = "mission"
$ws2.Range("C7").Value = "pratique délibérée"
$ws2.Range("D7").Value = 44260
$ws2.Range("E7").Value = 44261
$ws2.Range("F7").Value = 6
$ws2.Range("H7").Value = 4

# --- Row 8: intégration du titre du block hosting ---
$ws2.Range("A8").Value = "intégration du titre du block hosting"
$ws2.Range("B8").Value = "mission"
$ws2.Range("C8").Value = "pratique délibérée"
$ws2.Range("D8").Value = 44261
$ws2.Range("E8").Value = 44261
$ws2.Range("F8").Value = 1

# --- Row 9: intégration des cartes du block hosting ---
$ws2.Range("A9").Value = "intégration des cartes du block hosting"
$ws2.Range("B9").Value = "mission"
$ws2.Range("C9").Value = "pratique délibérée"
$ws2.Range("D9").Value = 44262
$ws2.Range("E9").Value = 44263
$ws2.Range("F9").Value = 6

# --- Row 10: flex froggy ---
$ws2.Range("A10").Value = "flex froggy"
$ws2.Range("B10").Value = "openclassrooms"
$ws2.Range("C10").Value = "activité complémentaire"
$ws2.Range("D10").Value = 44259
$ws2.Range("E10").Value = 44259
$ws2.Range("F10").Value = 1

# Widen column A on "étape 2" so the longer task descriptions fit.
$ws2.Columns.Item(1).ColumnWidth = 44.83

# Update the remembered selection on "étape 1" (user ended up on G15)
# without leaving that sheet active/selected in the saved workbook.
$ws1.Range("G15").Select()
$ws2.Activate()

# Update the remembered selection on "étape 2" (user ended up on H8:H9).
# "étape 2" stays the active/selected sheet, matching the saved file.
$ws2.Range("H8:H9").Select()
